$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Industry") - everything from C onward shifts right.
$ws.Columns("C:C").Insert()

# Header
$ws.Range("C1").Value = "Industry"

# Industry values for rows 2-19 (matches shifted data, column C rows 2..19)
$industries = @(
    "Pharmaceuticals & Biotechnology",
    "Food Products",
    "Pharmaceuticals & Biotechnology",
    "Pharmaceuticals & Biotechnology",
    "Healthcare Services",
    "Pharmaceuticals & Biotechnology",
    "Pharmaceuticals & Biotechnology",
    "Pharmaceuticals & Biotechnology",
    "Metals & Minerals Trading",
    "Pharmaceuticals & Biotechnology",
    "Retailing",
    "Pharmaceuticals & Biotechnology",
    "Pharmaceuticals & Biotechnology",
    "Pharmaceuticals & Biotechnology",
    "Pharmaceuticals & Biotechnology",
    "Realty",
    "Auto Components",
    "Pharmaceuticals & Biotechnology"
)

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $industries[$i]
}
